$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.329.42'
$ws.Range("E2").Value = '  -1.43%  '
$ws.Range("D3").Value = '2.918.60'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '364.49'
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.70'
$ws.Range("E6").Value = '  -5.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.540'
$ws.Range("E7").Value = '  -4.95%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.592'
$ws.Range("E9").Value = '  -6.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.01'
$ws.Range("E10").Value = '  -4.98%  '
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0835'
$ws.Range("E12").Value = '  -3.89%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.46'
$ws.Range("E13").Value = '  -5.46%  '
$ws.Range("D14").Value = '3.373.10'
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.36'
$ws.Range("E15").Value = '  -5.41%  '
$ws.Range("D16").Value = '2.918.32'
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.954'
$ws.Range("E17").Value = '  -3.47%  '
$ws.Range("D18").Value = '51.212.29'
$ws.Range("E18").Value = '  -1.65%  '
$ws.Range("E19").Value = '  -3.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.26'
$ws.Range("E20").Value = '  -3.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.06'
$ws.Range("E21").Value = '  -5.96%  '
$ws.Range("D22").Value = '0.0₃0949'
$ws.Range("E22").Value = '  -3.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.26'
$ws.Range("E23").Value = '  -3.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '260.38'
$ws.Range("E24").Value = '  -3.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.70'
$ws.Range("E25").Value = '  -4.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.32'
$ws.Range("E26").Value = '  +3.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.175'
$ws.Range("E27").Value = '  -4.97%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.95'
$ws.Range("E29").Value = '  -3.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.22'
$ws.Range("E30").Value = '  -6.03%  '
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("B32").Value = 'RenderToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.16'
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.96'
$ws.Range("E33").Value = '  -5.18%  '
$ws.Range("E34").Value = '  -3.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '35.14'
$ws.Range("E35").Value = '  -6.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '50.80'
$ws.Range("E36").Value = '  -2.60%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0424'
$ws.Range("E38").Value = '  -3.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.82'
$ws.Range("E39").Value = '  +3.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.15'
$ws.Range("E40").Value = '  -1.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.00'
$ws.Range("E41").Value = '  -6.88%  '
$ws.Range("E42").Value = '  -6.44%  '
$ws.Range("E43").Value = '  -5.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.46'
$ws.Range("E44").Value = '  -2.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '118.87'
$ws.Range("E45").Value = '  -0.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.12'
$ws.Range("E46").Value = '  -2.78%  '
$ws.Range("D47").Value = '2.061.94'
$ws.Range("E47").Value = '  -3.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.20'
$ws.Range("E48").Value = '  -7.42%  '
$ws.Range("E49").Value = '  -8.43%  '
$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.237'
$ws.Range("E50").Value = '  -6.83%  '
$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0315'
$ws.Range("E51").Value = '  -7.94%  '
